$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.20"
$ws.Range("E2").Value = "'-0.58%"
$ws.Range("D3").Value = "'26.25"
$ws.Range("E3").Value = "'3.04%"
$ws.Range("D4").Value = "'5.128"
$ws.Range("E4").Value = "'0.74%"
$ws.Range("D5").Value = "'0.05593"
$ws.Range("E5").Value = "'0.39%"
$ws.Range("D6").Value = "'6.475"
$ws.Range("E6").Value = "'-0.33%"
$ws.Range("D7").Value = "'0.8215"
$ws.Range("E7").Value = "'0.30%"
$ws.Range("E8").Value = "'-0.97%"
$ws.Range("E9").Value = "'-1.01%"
$ws.Range("D10").Value = "'0.06995"
$ws.Range("E10").Value = "'0.65%"
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = "'0.03119"
$ws.Range("E11").Value = "'-2.17%"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.02888"
$ws.Range("E12").Value = "'0.89%"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.09383"
$ws.Range("E13").Value = "'0.01%"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001513"
$ws.Range("E14").Value = "'0.20%"
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = "'0.0005991"
$ws.Range("E15").Value = "'-93.88%"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = "'0.006251"
$ws.Range("E16").Value = "'2.37%"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.653"
$ws.Range("E17").Value = "'4.40%"
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = "'3.033"
$ws.Range("E18").Value = "'0.47%"
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = "'2.183"
$ws.Range("E19").Value = "'4.38%"
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = "'0.3112"
$ws.Range("E20").Value = "'-2.12%"
$ws.Range("D21").Value = "'0.1300"
$ws.Range("E21").Value = "'-2.18%"
$ws.Range("D22").Value = "'3.741"
$ws.Range("E22").Value = "'-0.73%"
$ws.Range("D23").Value = "'0.04635"
$ws.Range("E23").Value = "'-1.63%"
$ws.Range("E24").Value = "'-0.07%"
$ws.Range("D25").Value = "'0.001245"
$ws.Range("E25").Value = "'-0.16%"
$ws.Range("E26").Value = "'-2.87%"
$ws.Range("D27").Value = "'0.00009598"
$ws.Range("E27").Value = "'-1.10%"
$ws.Range("D28").Value = "'0.0001394"
$ws.Range("E28").Value = "'0.21%"
$ws.Range("D40").Value = "'0.03642"
$ws.Range("E40").Value = "'-0.57%"
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = "'0.006179"
$ws.Range("E41").Value = "'-0.56%"
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1051"
$ws.Range("E42").Value = "'-0.08%"
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = "'0.002400"
$ws.Range("E43").Value = "'19.92%"
$ws.Range("D44").Value = "'0.008855"
$ws.Range("E44").Value = "'6.45%"
$ws.Range("D45").Value = "'0.00005344"
$ws.Range("E45").Value = "'0.65%"
$ws.Range("E46").Value = "'-0.07%"
$ws.Range("E47").Value = "'8.22%"
$ws.Range("D48").Value = "'0.002297"
$ws.Range("E48").Value = "'8.20%"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("E50").Value = "'-0.07%"
